# Insert two new rows before row 547, shifting existing rows 547:649 down to 549:651.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("547:548").Insert()

# Populate new row 547
$ws.Range("A547").Value = 4
$ws.Range("B547").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C547").Value = "Los Lagos"
$ws.Range("D547").Value = 45209
$ws.Range("E547").Value = 10
$ws.Range("F547").Value = 100114013
$ws.Range("G547").Value = "Zanahoria"
$ws.Range("H547").Value = "Sin especificar"
$ws.Range("I547").Value = "Primera"
$ws.Range("J547").Value = 350
$ws.Range("K547").Value = 8500
$ws.Range("L547").Value = 8500
$ws.Range("M547").Value = 8500
$ws.Range("N547").Value = "$/saco 20 kilos"
$ws.Range("O547").Value = "Provincia de Llanquihue"
$ws.Range("P547").Value = 425
$ws.Range("Q547").Value = 20
$ws.Range("R547").Value = "Hortaliza"

# Populate new row 548
$ws.Range("A548").Value = 4
$ws.Range("B548").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C548").Value = "Los Lagos"
$ws.Range("D548").Value = 45209
$ws.Range("E548").Value = 10
$ws.Range("F548").Value = 100114013
$ws.Range("G548").Value = "Zanahoria"
$ws.Range("H548").Value = "Sin especificar"
$ws.Range("I548").Value = "Primera"
$ws.Range("J548").Value = 350
$ws.Range("K548").Value = 9500
$ws.Range("L548").Value = 9500
$ws.Range("M548").Value = 9500
$ws.Range("N548").Value = "$/saco 20 kilos"
$ws.Range("O548").Value = "Región de La Araucanía"
$ws.Range("P548").Value = 475
$ws.Range("Q548").Value = 20
$ws.Range("R548").Value = "Hortaliza"
